# Update class diagram for Model, Storage and UI in dev guide
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)

    # "UI" group container (roundRect "Rectangle 65") - widen it and shift it
    # left so the relocated boxes inside still fit.
    if ($shape.Name -eq "Rectangle 65" -and $shape.TextFrame.TextRange.Text -eq "UI") {
        $shape.Left = 762001 / 12700
        $shape.Width = 5372548 / 12700
    }

    if ($shape.HasTextFrame) {
        $text = $shape.TextFrame.TextRange.Text
        if ($text -eq "PersonListPanel") {
            $shape.TextFrame.TextRange.Text = "CardListPanel"
        }
        elseif ($text -eq "PersonCard") {
            $shape.TextFrame.TextRange.Text = "CardDisplay"
        }
    }
}
